# Regenerate save_data to use K instead of Strike#: update column G ("K")
# values on the active sheet to the recalculated K counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K value (column G)
$kValues = @{
    2  = 0
    3  = 3
    4  = 1
    5  = 1
    6  = 1
    7  = 2
    8  = 2
    9  = 2
    11 = 0
    13 = 1
    14 = 2
    15 = 0
    16 = 2
    17 = 0
    18 = 1
    19 = 1
    20 = 0
    21 = 0
    22 = 0
    23 = 1
    24 = 1
    25 = 1
    26 = 0
    27 = 0
    28 = 1
    29 = 0
    30 = 1
    31 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
